$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '69.423.43'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.422.15'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.29'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.80'
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.414.64'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.198'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '48.79'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '689.51'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = '3.971.06'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = '69.465.74'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '3.424.03'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '100.68'
$ws.Range('E25').Value = '  -3.63%  '
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.58'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.40'
$ws.Range('E29').Value = '  -3.23%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.99'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '570.90'
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.68'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '58.29'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  -3.14%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '3.587.41'
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.139'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '34.86'
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  +2.63%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.26'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('E44').Value = '  -2.56%  '
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.65'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '131.88'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.63'
$ws.Range('E51').Value = '  +1.05%  '
